$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: 150001g1 / Dương Thị Hiệp / 1977 / father=150001 / mother=150001v1
$ws.Range("A22").Value = "150001g1"
$ws.Range("B22").Value = "Dương Thị Hiệp"
$ws.Range("C22").Value = 1977
$ws.Range("F22").Value = 150001
$ws.Range("G22").Value = "150001v1"

# Row 23: 150001g2 / Dương Thị Hường / 1991 / father=150001 / mother=150001v1
$ws.Range("A23").Value = "150001g2"
$ws.Range("B23").Value = "Dương Thị Hường"
$ws.Range("C23").Value = 1991
$ws.Range("F23").Value = 150001
$ws.Range("G23").Value = "150001v1"

# Row 24: 150002g1 / Dương Thị Thu / 1978 / father=150002 / mother=150002v1
$ws.Range("A24").Value = "150002g1"
$ws.Range("B24").Value = "Dương Thị Thu"
$ws.Range("C24").Value = 1978
$ws.Range("F24").Value = 150002
$ws.Range("G24").Value = "150002v1"

# Match the style (center aligned, like other A/F/G ID columns) for the new cells
$ws.Range("A22:A24").HorizontalAlignment = -4108
$ws.Range("F22:F24").HorizontalAlignment = -4108
$ws.Range("G22:G24").HorizontalAlignment = -4108

# Update selection to mirror the after-state
$ws.Range("E25").Select()
